# "updated entries till 13th july morning"
# Appends two new logBook entries (rows 53 & 54) for 13-Jul-2022, shifting
# the running "Total Hours" sum, and moves the sheet's active selection
# down to where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 53: 13-Jul-2022, 05:00-06:00, Code -------------------------------
# Copy formatting down from row 52 (the prior data row) so the new row
# picks up the same number formats / alignment styles used throughout the
# log (date, time x2, duration formula, category, wrapped description).
$ws.Range("A52:G52").Copy()
$ws.Range("A53:G53").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A53").Value2 = 52
$ws.Range("B53").Value2 = 44755
$ws.Range("C53").Value2 = 0.20833333333333334
$ws.Range("D53").Value2 = 0.25
$ws.Range("E53").Formula = "=D53-C53"
$ws.Range("F53").Value2 = "Code"
$ws.Range("G53").Value2 = "1. Literature survey on training segformer`n2. Understand Poly LR scheduler with 12ep training`n3. optimizer used - AdamW"

# This row's wrapped description spans 3 visual lines - size the row to
# match (45pt = 3 x 15pt default row height), same as the rest of the sheet.
$ws.Rows("53:53").RowHeight = 45

# --- Row 54: 13-Jul-2022, 07:45-08:20, Code -------------------------------
# Row 51 has no extra row height (single-line description), matching what
# row 54 needs, so use it as the formatting template.
$ws.Range("A51:G51").Copy()
$ws.Range("A54:G54").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A54").Value2 = 53
$ws.Range("B54").Value2 = 44755
$ws.Range("C54").Value2 = 0.32291666666666669
$ws.Range("D54").Value2 = 0.34722222222222227
$ws.Range("E54").Formula = "=D54-C54"
$ws.Range("F54").Value2 = "Code"
$ws.Range("G54").Value2 = "1. Poly LR scheduler, AdamW optimizer with 12ep training"

$excel.CutCopyMode = 0

# The totals formula in E57 (SUM(E2:E53)) auto-extends to include the new
# rows once row 53 is inserted before it, recalculating the grand total.

# Move the active selection / view to where the user left off editing.
[void]$ws.Range("G55").Select()
